# Auto-generated edit script applying the commit diff to cryptos.xlsx
# Updates Price (D) and Volume(1h) (E) columns for rows 2-45,
# and shifts the coin list in rows 46-51 (EnergySwap dropped, ordi added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.176.27'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '2.251.77'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.91'
$ws.Range("E5").Value = '  -4.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.55'
$ws.Range("E6").Value = '  -2.30%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -3.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.56'
$ws.Range("E10").Value = '  -3.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0821'
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.31'
$ws.Range("E12").Value = '  -5.40%  '
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").Value = '2.594.51'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '2.287.86'
$ws.Range("E15").Value = '  +2.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.837'
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.79'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '43.995.42'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.86'
$ws.Range("E19").Value = '  -4.74%  '
$ws.Range("D20").Value = '0.0₃0974'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.33'
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.41'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("E24").Value = '  -7.03%  '
$ws.Range("E25").Value = '  -7.72%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.62'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.14'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.54'
$ws.Range("E32").Value = '  -2.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.50'
$ws.Range("E33").Value = '  +13.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0820'
$ws.Range("E34").Value = '  -3.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.66'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("E37").Value = '  -4.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  -3.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.44'
$ws.Range("E39").Value = '  -0.35%  '
$ws.Range("E40").Value = '  -8.66%  '
$ws.Range("E41").Value = '  -3.24%  '
$ws.Range("E42").Value = '  -10.15%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '1.763.25'
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '86.96'
$ws.Range("E45").Value = '  +5.70%  '
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.16'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.193'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.29'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.25'
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.74'
$ws.Range("E50").Value = '  -4.72%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.37'
$ws.Range("E51").Value = '  -6.89%  '
